# Add team record (Wins / Losses / Ties) columns to the right of the
# existing data, matching the formatting of the last existing header
# column (AC1) for the new header cells AD1:AF1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font, border, centered alignment)
# from the last header cell onto the three new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-50) gets the team's 1999 record: 94-68-0.
$ws.Range("AD2:AD50").Value = 94
$ws.Range("AE2:AE50").Value = 68
$ws.Range("AF2:AF50").Value = 0
